$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column header: I1 "Commission" (copy the bold/bordered header format from H1) ---
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("I1").Value = "Commission"

function Set-TextValue($addr, $text) {
    # Force a numeric-looking string (e.g. "35.0", "19.579%") to be stored as
    # literal text instead of being auto-converted to a number/percentage by
    # Excel's input parser, then restore the cell to the default "Normal"
    # style so no stray number-format is left behind on the cell.
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

# --- Row 2: EREGL.IS ---
$ws.Range("A2").Value = "EREGL.IS"
Set-TextValue "B2" "35.0"
Set-TextValue "C2" "32.163"
$ws.Range("D2").Value = 1125.7
Set-TextValue "E2" "38.46"
$ws.Range("F2").Value = 1346.1
$ws.Range("G2").Value = 220.4
Set-TextValue "H2" "19.579%"
$ws.Range("I2").Value = 2.35

# --- Row 3: SISE.IS (new row) ---
$ws.Range("A3").Value = "SISE.IS"
Set-TextValue "B3" "30.0"
Set-TextValue "C3" "34.107"
$ws.Range("D3").Value = 1023.2
Set-TextValue "E3" "49.94"
$ws.Range("F3").Value = 1498.2
$ws.Range("G3").Value = 475
Set-TextValue "H3" "46.423%"
$ws.Range("I3").Value = 2.14

# --- Row 4: THYAO.IS (new row) ---
$ws.Range("A4").Value = "THYAO.IS"
Set-TextValue "B4" "3.0"
Set-TextValue "C4" "107.7"
$ws.Range("D4").Value = 323.1
Set-TextValue "E4" "235.0"
$ws.Range("F4").Value = 705
$ws.Range("G4").Value = 381.9
Set-TextValue "H4" "118.199%"
$ws.Range("I4").Value = 0.68

# --- Row 5: TUKAS.IS (new row) ---
$ws.Range("A5").Value = "TUKAS.IS"
$ws.Range("B5").Value = 25
$ws.Range("C5").Value = 20.32
$ws.Range("D5").Value = 508
$ws.Range("E5").Value = 10.08
$ws.Range("F5").Value = 252
$ws.Range("G5").Value = -256
Set-TextValue "H5" "-50.394%"
$ws.Range("I5").Value = 1.06

# --- Row 6: Totals row (new row) ---
$ws.Range("A6").Value = "-"
$ws.Range("B6").Value = "-"
$ws.Range("C6").Value = "-"
$ws.Range("D6").Value = 2980
$ws.Range("E6").Value = "-"
$ws.Range("F6").Value = 3801.3
$ws.Range("G6").Value = 821.3
Set-TextValue "H6" "27.56%"
$ws.Range("I6").Value = 6.23
